$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1
$ws.Range("I1:J1").Borders.Weight = 2

$data = @(
    5,6
    8,9
    8,8
    9,9
    7,8
    8,8
    5,6
    8,9
    9,9
    5,6
    6,8
    5,5
    8,8
    5,6
    8,8
    8,8
    7,7
    8,9
    7,8
    10,10
    9,9
    5,6
    8,9
    7,8
    7,8
    9,9
    4,5
    8,8
    7,8
    8,8
    7,8
    8,8
    6,6
    6,6
    3,3
)

for ($i = 0; $i -lt $data.Length / 2; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i * 2]
    $ws.Cells.Item($row, 10).Value = $data[$i * 2 + 1]
}
